# Fix instructor view: rename the GA (grading area) header columns on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header row (row 1) GA column labels D1:H1
$ws.Range("D1").Value = "GA 3.1"
$ws.Range("E1").Value = "GA 8.1"
$ws.Range("F1").Value = "GA 10.3"
$ws.Range("G1").Value = "GA 5.2"
$ws.Range("H1").Value = "GA 4.1"

# Move/restore the active selection to H1 (matches the saved view state)
$ws.Range("H1").Select()
